$d = $word.ActiveDocument
$d.Content.Find.Execute("LUCHAO QI", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Luchao Qi", 2)
